# Update column G ("K") values on Sheet1 with newly computed strikeout
# counts (regenerated from pitch-level data: save_data now uses K instead
# of Strike# for this column, and std/mean + s_vals were recalculated
# upstream). Only the G column values for rows 2-15 change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @{
    2  = 2
    3  = 0
    4  = 3
    5  = 0
    6  = 2
    7  = 3
    8  = 7
    9  = 5
    10 = 1
    11 = 2
    12 = 4
    13 = 3
    14 = 1
    15 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
